$wb = $excel.ActiveWorkbook

# Sheet "展览": F3 275 -> 277, F4 2705 -> 2717
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 277
$ws1.Range("F4").Value = 2717

# Sheet "演出": F3 18 -> 19
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 19

# Sheet "全部类型": F3 18 -> 19, F5 275 -> 277, F6 2705 -> 2717
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 19
$ws4.Range("F5").Value = 277
$ws4.Range("F6").Value = 2717
